$wb = $excel.ActiveWorkbook

# --- 1. Update status text "Ready for handoff" -> "In Translation" ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws1.Range("E3").Value = "In Translation"
$ws1.Range("F3").Value = "In Translation"

$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2").Value = "In Translation"
$ws2.Range("C3").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2").Value = "In Translation"
$ws3.Range("C3").Value = "In Translation"

# --- 2. Narrow status columns (previously ~17.22 chars wide, now ~13.41 chars wide) ---
# Target stored width is 13.4101845877511 characters. This runtime quantizes the
# ColumnWidth setter to 1/6-character steps, so 12.5 is the input that lands on the
# closest achievable stored width (13.3333...).
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws3.Columns.Item(3).ColumnWidth = 12.5
